$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-216). The sheet was refreshed and every one of these dates
# moved forward by a single day: 45188 (2023-09-19) -> 45189 (2023-09-20).
$lastRow = 216
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
